$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("C1").Value = "cp"

# Update city names in column B
$ws.Range("B2").Value = "miguel hidalgo "
$ws.Range("B3").Value = "xochimilco "

# Fix state value in A3 (was 2, should be 1)
$ws.Range("A3").Value = 1

# Add new postal code (cp) column values
$ws.Range("C2").Value = 11000
$ws.Range("C3").Value = 16000

# Update the selected/active cell to match target state
$ws.Range("C2").Select()
